$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.780.86'
$ws.Range('E2').Value = '  -0.42%  '

$ws.Range('D3').Value = '1.683.47'
$ws.Range('E3').Value = '  -1.41%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9997'
$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.87'
$ws.Range('E5').Value = '  -0.97%  '

$ws.Range('E6').Value = '  +0.25%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3921'
$ws.Range('E7').Value = '  -1.99%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3970'
$ws.Range('E8').Value = '  -2.80%  '

$ws.Range('B9').Value = 'OKB'
$ws.Range('C9').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '51.99'
$ws.Range('E9').Value = '  -3.09%  '

$ws.Range('B10').Value = 'BinanceUSD'
$ws.Range('C10').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9993'
$ws.Range('E10').Value = '  -0.15%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.409'
$ws.Range('E11').Value = '  -5.23%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08674'
$ws.Range('E12').Value = '  -2.05%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '25.37'
$ws.Range('E13').Value = '  -3.73%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.347'
$ws.Range('E14').Value = '  -2.11%  '

$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001320'
$ws.Range('E15').Value = '  -3.11%  '

$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.777'
$ws.Range('E16').Value = '  -4.68%  '

$ws.Range('D17').Value = '1.759.61'
$ws.Range('E17').Value = '  +3.44%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '93.84'
$ws.Range('E18').Value = '  -3.33%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.07087'
$ws.Range('E19').Value = '  -1.39%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '20.32'
$ws.Range('E20').Value = '  -4.66%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.084'
$ws.Range('E21').Value = '  -2.98%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.002'
$ws.Range('E22').Value = '  +0.45%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.97'
$ws.Range('E23').Value = '  -3.03%  '

$ws.Range('D24').Value = '24.763.48'
$ws.Range('E24').Value = '  -0.42%  '

$ws.Range('E25').Value = '  +1.04%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.839'
$ws.Range('E26').Value = '  -3.80%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '23.56'
$ws.Range('E27').Value = '  +0.94%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '162.75'
$ws.Range('E28').Value = '  -2.49%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.839'
$ws.Range('E29').Value = '  -7.32%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '146.97'
$ws.Range('E30').Value = '  +0.45%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.917'
$ws.Range('E31').Value = '  -6.07%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.387'
$ws.Range('E32').Value = '  +6.63%  '

$ws.Range('D33').Value = '1.866.02'
$ws.Range('E33').Value = '  -1.56%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08448'
$ws.Range('E34').Value = '  -4.61%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.03060'
$ws.Range('E35').Value = '  -4.61%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.984'
$ws.Range('E36').Value = '  -3.16%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.004'
$ws.Range('E37').Value = '  -3.65%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2813'
$ws.Range('E38').Value = '  -2.71%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.09467'
$ws.Range('E39').Value = '  +1.14%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '10.65'
$ws.Range('E40').Value = '  -2.25%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.519'
$ws.Range('E41').Value = '  +3.35%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.7965'
$ws.Range('E42').Value = '  -7.06%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '13.59'
$ws.Range('E43').Value = '  -4.56%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.65'
$ws.Range('E44').Value = '  -5.09%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.7168'
$ws.Range('E45').Value = '  -4.29%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.575'
$ws.Range('E46').Value = '  -5.49%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.187'
$ws.Range('E47').Value = '  -1.37%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.08699'
$ws.Range('E48').Value = '  +4.09%  '

$ws.Range('B49').Value = 'Flow'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.353'
$ws.Range('E49').Value = '  -3.62%  '

$ws.Range('B50').Value = 'Frax'
$ws.Range('C50').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.001'
$ws.Range('E50').Value = '  +0.30%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '138.60'
$ws.Range('E51').Value = '  -2.24%  '
